$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.450.95"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.571.64"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'291.64"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("D8").Value = "'49.83"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.3405"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "'0.07551"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'6.039"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "'6.962"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "1.589.82"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "'0.00001122"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "'91.22"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "'0.06759"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'6.300"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "'16.29"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "'12.16"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "22.443.66"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'2.340"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").Value = "'2.651"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "'148.83"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "'5.035"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "'125.84"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "1.754.89"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "'1.068"
$ws.Range("E32").Value = "  +8.40%  "
$ws.Range("D33").Value = "'6.205"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "'2.001"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'9.826"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").Value = "'0.08381"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").Value = "'0.02484"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.350"
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2302"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'0.06531"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'5.462"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "'11.30"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "'0.6229"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'14.04"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "'3.806"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5813"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'130.22"
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").Value = "'1.226"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("D51").Value = "'0.07319"
$ws.Range("E51").Value = "  -0.22%  "
